# Kammari_LabExam03Grading.xlsx - grading pass ("kalyankar to pusapati done")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Customer Class section (rows 3-6) ---
$ws.Range("F3").Value = "(-1) for missing author notation"
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Product Class section (rows 10-14) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the active selection to E15 (previously scrolled to A24 / C44 selected)
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("E15").Select()
